# Updating readme with more samples + DNode
#
# On the "Who you gonna call?" slide, the content placeholder lists a
# bunch of npm package names separated by ", ". In the authored edit the
# run boundary between "...pause-stream" and the following ", "/"emit-stream"
# runs is redrawn: the trailing ", " is folded into the "...pause-stream"
# run (giving "...pause-stream, "), and "emit-stream" becomes its own run
# immediately after (picking up a fresh dirty="0" the way PowerPoint marks
# text it just reflowed).

$p = $ppt.ActivePresentation

# Find the slide that holds the "pause-stream" list (content placeholder).
$slide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    foreach ($shp in $candidate.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text.IndexOf("pause-stream") -ge 0) {
            $slide = $candidate
            break
        }
    }
    if ($slide -ne $null) { break }
}

$shape = $null
foreach ($shp in $slide.Shapes) {
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text.IndexOf("pause-stream") -ge 0) {
        $shape = $shp
        break
    }
}

$tr = $shape.TextFrame.TextRange
$full = $tr.Text

# --- locate the two text segments involved in the edit -------------------
$tailMarker  = ", shoe, pause-stream"          # end of the run that keeps growing
$commaMarker = ", "                             # the lone ", " run right after it
$nextWord    = "emit-stream"                    # the run that should pick up dirty="0"

$tailStart = $full.IndexOf($tailMarker) + 1     # 1-based start for Characters()
$tailLen   = $tailMarker.Length

$commaStart = $tailStart + $tailLen
$commaLen   = $commaMarker.Length

$wordStart = $commaStart + $commaLen
$wordLen   = $nextWord.Length

# --- step 1: merge the ", " run into the "...pause-stream" run -----------
$tailRun = $tr.Characters($tailStart, $tailLen)
$tailRun.Text = $tailMarker + $commaMarker      # ", shoe, pause-stream, "

$commaRun = $tr.Characters($wordStart - $commaLen, $commaLen)
$commaRun.Text = ""                              # drop the now-empty run

# --- step 2: give "emit-stream" its own run stamped dirty="0" ------------
# Grab the trailing space of the merged run together with the following
# word; assigning new text to a multi-run selection collapses it into a
# single run that inherits the formatting (incl. dirty="0") of the first
# run touched.
$newTailLen = $tailLen + $commaLen               # length of ", shoe, pause-stream, "
$boundaryPos = $tailStart + $newTailLen - 1      # 1-based position of the trailing space
$combo = $tr.Characters($boundaryPos, 1 + $wordLen)
$combo.Text = " " + $nextWord

# Split the leading space back off so "emit-stream" stands alone.
$lead = $tr.Characters($boundaryPos, 1)
$lead.Text = ""

# Restore the trailing space on the "...pause-stream, " run (this is a
# same-run, in-place edit so it does not disturb the neighbouring run).
$tailRunFixed = $tr.Characters($tailStart, $newTailLen - 1)
$tailRunFixed.Text = $tailMarker + $commaMarker

$tr.Text
